$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rows 8 and 9 had their match data (columns F:V) swapped -----------
# (columns A:E — Indice/pais/torneio/temporada/data_partida — stay put)
$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

$row8vals = @{}
$row9vals = @{}
foreach ($c in $cols) {
  $row8vals[$c] = $ws.Range($c + "8").Value2
  $row9vals[$c] = $ws.Range($c + "9").Value2
}
foreach ($c in $cols) {
  $ws.Range($c + "8").Value = $row9vals[$c]
  $ws.Range($c + "9").Value = $row8vals[$c]
}

# --- 2) Append a new row 44 with a further match -----------------------
# Copy formatting from row 43 (the previous last row) for the styled cells
# (A = bold/bordered index style, E = date-time number format), then set values.
$ws.Range("A43").Copy()
$ws.Range("A44").PasteSpecial(-4122)
$ws.Range("E43").Copy()
$ws.Range("E44").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A44").Value = 43
$ws.Range("B44").Value = "moldova"
$ws.Range("C44").Value = "super-liga"
$ws.Range("D44").Value = "2023-2024"
$ws.Range("E44").Value = 45235.58333333334
$ws.Range("F44").Value = "Balti"
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = "Zimbru Chisinau"
$ws.Range("I44").Value = 1
$ws.Range("J44").Value = 2.39
$ws.Range("K44").Value = "04/11/2023 02:13"
$ws.Range("L44").Value = 3.27
$ws.Range("M44").Value = "05/11/2023 13:31"
$ws.Range("N44").Value = 2.9
$ws.Range("O44").Value = "04/11/2023 02:13"
$ws.Range("P44").Value = 3.12
$ws.Range("Q44").Value = "05/11/2023 13:31"
$ws.Range("R44").Value = 2.61
$ws.Range("S44").Value = "04/11/2023 02:13"
$ws.Range("T44").Value = 2.09
$ws.Range("U44").Value = "05/11/2023 13:31"
$ws.Range("V44").Value = "https://www.betexplorer.com/football/moldova/super-liga/csf-balti-zimbru-chisinau/CGgaQ4lq/"
